$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are a data engineer at a multinational company that collects data from multiple sources like IoT devices, social media, and internal databases. The data varies in format from structured to unstructured types, and you need to set up an ETL architecture for it.Which ETL architecture should you use?",
        "ques_type": 2,
        "options": [
            "Single tier",
            "Two-tier",
            "On-premises",
            "Off-premises"
        ],
        "score": "Two-tier"
    },
    {
        "title": "You are a data engineer managing an ETL process for an e-commerce company. Recently, the transformation stage of the ETL process has been taking a significantly longer amount of time due to increased data volume and complexity. You have been tasked with improving the performance of this step. What should you do?",
        "ques_type": 2,
        "options": [
            "Simplify or optimize the transformation logic.",
            "Increase the network bandwidth.",
            "Upgrade the storage capacity of the data warehouse.",
            "Reschedule ETL processes to avoid competition."
        ],
        "score": "Simplify or optimize the transformation logic."
    },
    {
        "title": "The performance of your organization's data pipeline has recently had significant degradation, affecting the overall data processing speed. After investigating, you suspect the problem might lie in the data transformation phase. What should you do?",
        "ques_type": 2,
        "options": [
            "Increase the processing power of the system.",
            "Profile the transformation operations.",
            "Redefine the transformation operations.",
            "Delete old records."
        ],
        "score": "Profile the transformation operations."
    },
    {
        "title": "You have been tasked with designing a failure recovery strategy for your company\u2019s ETL pipeline. The pipeline handles multiple large datasets and performs a series of complex transformations. Which of the following is the best general solution?",
        "ques_type": 2,
        "options": [
            "Update the ETL software regularly.",
            "Implement a checkpoint mechanism.",
            "Use high-performance hardware.",
            "Use additional computational resources during peak hours."
        ],
        "score": "Implement a checkpoint mechanism."
    }
]
'@

# The here-string adds exactly one trailing newline after the closing
# bracket; strip it so the cell content matches the source text exactly.
$newText = $newText.TrimEnd("`r", "`n")

# Row 1 previously held a plain "0" with a bold/bordered/centered style, and
# the real question text lived in A2 as a shared string. Clear both cells
# and their formatting, then write the (reformatted) text into A1 alone.
$ws.Range("A2").ClearContents()
$ws.Range("A1:A2").ClearFormats()
$ws.Range("A1").Value = $newText
